# Common: Prepared very first transaction stuff
# Adds a "cost" column (G) to the "atomizers" worksheet with per-row costs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("atomizers")

# Header - match the bold "Nadpis 2" heading style used by the rest of row 1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "cost"

# Row costs (rows 12 and 13 intentionally have no cost yet - still touched as
# blanks so their row spans get bumped, matching row 18 which is left
# completely untouched and keeps its original span)
$costs = @{
    2 = 50
    3 = 150
    4 = 250
    5 = 25
    6 = 75
    7 = 275
    8 = 10
    9 = 300
    10 = 275
    11 = 275
    12 = ""
    13 = ""
    14 = 75
    15 = 275
    16 = 75
    17 = 275
}

foreach ($row in 2..17) {
    $ws.Cells.Item($row, 7).Value = $costs[$row]
}

# Update selection to match the authored workbook state
$ws.Range("G17").Select()
